# Insert a new weekly data row for "Espinaca" (Vega Modelo de Temuco) at
# row 86, pushing the existing rows 86-114 down to 87-115.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(86).Insert()

# Populate the newly inserted row 86 with the new record.
$ws.Cells.Item(86, 1).Value = 10
$ws.Cells.Item(86, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(86, 3).Value = "La Araucanía"
$ws.Cells.Item(86, 4).Value = 44553
$ws.Cells.Item(86, 5).Value = 9
$ws.Cells.Item(86, 6).Value = 100112012
$ws.Cells.Item(86, 7).Value = "Espinaca"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 135
$ws.Cells.Item(86, 11).Value = 8000
$ws.Cells.Item(86, 12).Value = 9000
$ws.Cells.Item(86, 13).Value = 8481
$ws.Cells.Item(86, 14).Value = "$/docena de atados"
$ws.Cells.Item(86, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(86, 16).Value = 2827
$ws.Cells.Item(86, 17).Value = 3
$ws.Cells.Item(86, 18).Value = "Hortaliza"
